# Adds a "one-hot" breakdown of the final grade (column M, values 2-5) into
# four new columns R:U - one column per possible grade value - so R3:U3 hold
# the headers 2,3,4,5 and, for every student row (4-32), exactly one of
# R:U is 1 (the others 0) depending on which grade that student got.
# Also zooms out to 70%, moves the selection onto the new block, and paints
# it with the same 3-colour colour-scale conditional formatting already used
# on the other grade columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: the 4 possible grade values -----------------------------
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 5

# --- one-hot indicator formulas for every student row (4-32) ------------
# Entering the formula on the whole block at once makes Excel store it as a
# single shared-formula group, each cell comparing its own column header
# (row 3, absolute row) against that row's grade in column M (absolute col).
$ws.Range("R4:U32").Formula = "=IF(R`$3=`$M4,1,0)"

# --- view: zoom out and select the new block -----------------------------
$ws.Application.ActiveWindow.Zoom = 70
[void]$ws.Range("R4:U32").Select()

# --- conditional formatting: same 3-colour colour scale as columns C:M --
# (red -> yellow -> green, split at min / 50th percentile / max), added as
# a brand-new rule so it is bumped to first priority like Excel does when a
# new conditional-formatting rule is created through the UI.
$cf = $ws.Range("R4:U32").FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria.Item(1).Type = 1               # xlConditionValueLowestValue
$cf.ColorScaleCriteria.Item(2).Type = 5               # xlConditionValuePercentile
$cf.ColorScaleCriteria.Item(2).Value = 50
$cf.ColorScaleCriteria.Item(3).Type = 2               # xlConditionValueHighestValue
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 7039480   # RGB(248,105,107) = F8696B
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167   # RGB(255,235,132) = FFEB84
$cf.ColorScaleCriteria.Item(3).FormatColor.Color = 8109667   # RGB(99,190,123)  = 63BE7B
$cf.SetFirstPriority()
